# Fix config file paths - remove duplicate examples/ prefix
$wb = $excel.ActiveWorkbook

# --- Settings sheet: fix data_file / output_file paths ---
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("B8").Value = "sample_cbc_data.csv"
$settings.Range("B9").Value = "output/example_results.xlsx"

# --- Instructions sheet: clear stray empty inline-string cells on blank rows ---
$instructions = $wb.Worksheets.Item("Instructions")
$blankRows = @(2, 4, 10, 17, 23, 29)
foreach ($r in $blankRows) {
    $instructions.Cells.Item($r, 1).ClearContents()
}
